$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.372.63"
$ws.Range("E2").Value = "  -0.07%  "

# Row 3
$ws.Range("D3").Value = "2.067.55"
$ws.Range("E3").Value = "  +0.04%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.72"
$ws.Range("E5").Value = "  -0.36%  "

# Row 6
$ws.Range("E6").Value = "  +1.52%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.06"
$ws.Range("E8").Value = "  -2.61%  "

# Row 9
$ws.Range("E9").Value = "  +2.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0776"
$ws.Range("E10").Value = "  +1.72%  "

# Row 11
$ws.Range("E11").Value = "  +0.88%  "

# Row 12
$ws.Range("D12").Value = "2.372.61"
$ws.Range("E12").Value = "  +0.06%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.33"
$ws.Range("E13").Value = "  -1.78%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.63"
$ws.Range("E14").Value = "  -3.71%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.777"
$ws.Range("E15").Value = "  +0.00%  "

# Row 16
$ws.Range("E16").Value = "  +0.17%  "

# Row 17
$ws.Range("D17").Value = "2.068.25"
$ws.Range("E17").Value = "  -0.01%  "

# Row 18
$ws.Range("D18").Value = "37.326.62"
$ws.Range("E18").Value = "  -0.71%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.23"
$ws.Range("E19").Value = "  +2.40%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.52"
$ws.Range("E20").Value = "  -1.29%  "

# Row 21
$ws.Range("E21").Value = "  +0.40%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.30"
$ws.Range("E22").Value = "  -0.66%  "

# Row 23
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
$ws.Range("E24").Value = "  +1.63%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.39"
$ws.Range("E25").Value = "  -2.64%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.85"
$ws.Range("E26").Value = "  +1.39%  "

# Row 27
$ws.Range("E27").Value = "  -0.45%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.132"
$ws.Range("E28").Value = "  +3.93%  "

# Row 29
$ws.Range("E29").Value = "  -6.31%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.07"
$ws.Range("E30").Value = "  -0.98%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.117"
$ws.Range("E31").Value = "  -1.46%  "

# Row 32
$ws.Range("E32").Value = "  +0.52%  "

# Row 33
$ws.Range("E33").Value = "  -1.09%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.54"
$ws.Range("E34").Value = "  -0.55%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.46"
$ws.Range("E35").Value = "  -2.98%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.37"
$ws.Range("E36").Value = "  -0.71%  "

# Row 37
$ws.Range("E37").Value = "  -0.35%  "

# Row 38
$ws.Range("E38").Value = "  +0.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.61"
$ws.Range("E39").Value = "  -4.43%  "

# Row 40
$ws.Range("E40").Value = "  -0.34%  "

# Row 41
$ws.Range("D41").Value = "1.490.09"
$ws.Range("E41").Value = "  +2.21%  "

# Row 42
$ws.Range("E42").Value = "  -2.68%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.87"
$ws.Range("E43").Value = "  +1.45%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0213"
$ws.Range("E44").Value = "  +0.51%  "

# Row 45
$ws.Range("E45").Value = "  -0.77%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.22"
$ws.Range("E46").Value = "  -4.37%  "

# Row 47
$ws.Range("E47").Value = "  -0.16%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.21"
$ws.Range("E48").Value = "  -4.32%  "

# Row 49
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.22"
$ws.Range("E49").Value = "  -0.72%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.97"
$ws.Range("E50").Value = "  +0.76%  "

# Row 51
$ws.Range("D51").Value = "2.259.29"
$ws.Range("E51").Value = "  +0.03%  "
